# "add dynamic plot method in serial port" - extend the motor_data sample
# column A with additional simulated sensor readings (rows 324-625), then
# leave the view scrolled/selected near the newly appended data, matching
# where the user was working when the rows were added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three blocks of constant values below the existing data (which
# already ends at row 323 with value 5).
$ws.Range("A324:A360").Value = 5
$ws.Range("A361:A470").Value = 12
$ws.Range("A471:A625").Value = 0

# Scroll the view down toward the new rows and select the cell the user
# ended up on.
$excel.ActiveWindow.ScrollRow = 609
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G617").Select() | Out-Null
